# Add a new "localdb" command-category column to the hidden '#system' sheet,
# insert it into the alphabetical 'target' list, and wire up the named ranges
# (mirrors Nexial's convention: one column per command group, row1 = group
# name, rows2.. = command signatures; column A lists every group name and is
# used as the source for a dependent data-validation dropdown).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Insert a new column before N. This pushes the existing N:AC columns
#    (macro, mail, number, pdf, rdbms, redis, sms, sound, ssh, step, web,
#    webalert, webcookie, ws, ws.async, xml) one column to the right,
#    becoming O:AD.
# ---------------------------------------------------------------------
$ws.Columns("N:N").Insert()

# ---------------------------------------------------------------------
# 2) Populate the freed-up column N with the new 'localdb' command group.
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 14).Value = "localdb"
$ws.Cells.Item(2, 14).Value = "cloneTable(var,source,target)"
$ws.Cells.Item(3, 14).Value = "dropTables(var,tables)"
$ws.Cells.Item(4, 14).Value = "exportCSV(sql,output)"
$ws.Cells.Item(5, 14).Value = "importRecords(var,sourceDb,sql,table)"
$ws.Cells.Item(6, 14).Value = "purge(var)"
$ws.Cells.Item(7, 14).Value = "runSQLs(var,sqls)"

# ---------------------------------------------------------------------
# 3) Insert 'localdb' into column A's alphabetical list of group names
#    (the 'target' range), right before 'macro', shifting macro..xml
#    down by one row (A14:A29 -> A15:A30). Only column A moves; every
#    other column keeps its existing row alignment.
# ---------------------------------------------------------------------
for ($r = 29; $r -ge 14; $r--) {
    $v = $ws.Cells.Item($r, 1).Value2()
    $ws.Cells.Item($r + 1, 1).Value = $v
}
$ws.Cells.Item(14, 1).Value = "localdb"

# ---------------------------------------------------------------------
# 4) Fix up the defined names that pointed at the columns which just
#    shifted right, and register the new 'localdb' and 'target' ranges.
# ---------------------------------------------------------------------
function Set-NamedRange($name, $col, $rowStart, $rowEnd) {
    $n = $wb.Names.Item($name)
    $n.RefersTo = "='#system'!`$" + $col + "`$" + $rowStart + ":`$" + $col + "`$" + $rowEnd
}

Set-NamedRange "macro"     "O"  2 4
Set-NamedRange "mail"      "P"  2 2
Set-NamedRange "number"    "Q"  2 16
Set-NamedRange "pdf"       "R"  2 16
Set-NamedRange "rdbms"     "S"  2 7
Set-NamedRange "redis"     "T"  2 10
Set-NamedRange "sms"       "U"  2 2
Set-NamedRange "sound"     "V"  2 5
Set-NamedRange "ssh"       "W"  2 9
Set-NamedRange "step"      "X"  2 4
Set-NamedRange "web"       "Y"  2 127
Set-NamedRange "webalert"  "Z"  2 8
Set-NamedRange "webcookie" "AA" 2 8
Set-NamedRange "ws"        "AB" 2 17
Set-NamedRange "ws.async"  "AC" 2 8
Set-NamedRange "xml"       "AD" 2 21

Set-NamedRange "target" "A" 2 30

$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
